$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorder the "Periodo Mora" column (E16:E22) into ascending chronological
# order (2407 .. 2501) and carry the "Valor Mora" (F column) figure that
# belonged to period 2501 (50266) over to period 2407, while the former
# 2407 row now takes the "standard" 52000 value - i.e. swap F16 and F22.

$periods = @("2407", "2408", "2409", "2410", "2411", "2412", "2501")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
}

$ws.Range("F16").Value = 52000
$ws.Range("F22").Value = 50266
